# "Update evaluation for Stardog"
#
# Sheet1 used to hold per-reasoner timing/result rows (rows 2-12) plus a
# second summary block (rows 15-20) with SUM() formulas referencing those
# rows. The edit replaces all of that with a single, fresh data row for the
# new "Stardog" reasoner (row 5) and drops everything below it (rows 6-20,
# including the now-orphaned summary block), shrinking the sheet to A1:Q5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Drop rows 6:20 - the old extra reasoner rows (6-12) and the whole
# second summary table (15-20) that summed them.
$ws1.Rows("6:20").Delete()

# Row 5 becomes the "Stardog" results row.
$ws1.Range("A5").Value = "Stardog"
$ws1.Range("B5").Value = 1.156
$ws1.Range("C5").Value = 0.107
$ws1.Range("D5").Value = 2.489
$ws1.Range("E5").Value = 13.372
$ws1.Range("F5").Value = 40.219
$ws1.Range("G5").Value = 47.318

# Match the saved cursor position on each sheet.
$ws1.Range("G5").Select()

$ws2.Activate()
$ws2.Range("G14").Select()
